$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A11").Value = "20234211/24/23 -02:42:21"
$ws.Range("B11").Value = "ApplicationFrameHost.exe"
$ws.Range("C11").Value = 0
$ws.Range("D11").Value = 35.21875
$ws.Range("E11").Value = 0.2163152463620963
$ws.Range("F11").Value = 281
$ws.Range("G11").Value = 3

$ws.Range("A12").Value = "20234211/24/23 -02:42:58"
$ws.Range("B12").Value = "ApplicationFrameHost.exe"
$ws.Range("C12").Value = 0
$ws.Range("D12").Value = 35.21875
$ws.Range("E12").Value = 0.2163152463620963
$ws.Range("F12").Value = 558
$ws.Range("G12").Value = 3
